$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("G2").Value = 2.2
$ws.Range("I2").Value = 3.6
$ws.Range("L2").Value = 4.5
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 6
$ws.Range("AH2").Value = 17
$ws.Range("AK2").Value = 41

# Row 4 updates
$ws.Range("G4").Value = 2.3
$ws.Range("I4").Value = 3.4
$ws.Range("J4").Value = 3.25
$ws.Range("L4").Value = 4.33
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 1.62
$ws.Range("P4").Value = 2.2
$ws.Range("U4").Value = 2.38
$ws.Range("V4").Value = 1.53
$ws.Range("X4").Value = 9.5
$ws.Range("AH4").Value = 15
$ws.Range("AJ4").Value = 41
$ws.Range("AK4").Value = 41
$ws.Range("AN4").Value = 4
$ws.Range("AZ4").Value = 81
$ws.Range("BA4").Value = 151
